$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells before assigning numeric-looking strings,
# so Excel stores them as text instead of coercing to numbers.
$dCells = @("D2","D3","D4","D5","D6","D7","D9","D11","D14","D15","D16","D17","D18","D19","D21","D22","D23","D24","D25","D26","D27","D30","D31","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D46","D47","D48","D49","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row (order follows the source diff).
$ws.Range("D2").Value = '70.005.55'
$ws.Range("E2").Value = '  +1.26%  '

$ws.Range("D3").Value = '3.797.40'
$ws.Range("E3").Value = '  +0.35%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '663.23'
$ws.Range("E5").Value = '  +5.40%  '

$ws.Range("D6").Value = '166.54'
$ws.Range("E6").Value = '  +1.39%  '

$ws.Range("D7").Value = '3.796.32'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  +1.71%  '

$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").Value = '0.459'
$ws.Range("E11").Value = '  +1.81%  '

$ws.Range("E12").Value = '  +5.08%  '

$ws.Range("E13").Value = '  -2.41%  '

$ws.Range("D14").Value = '35.67'
$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").Value = '4.443.47'
$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16").Value = '3.801.87'
$ws.Range("E16").Value = '  +0.13%  '

$ws.Range("D17").Value = '69.990.36'
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").Value = '17.74'
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").Value = '7.17'
$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '474.58'
$ws.Range("E21").Value = '  +1.24%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '10.21'
$ws.Range("E22").Value = '  +5.73%  '

$ws.Range("D23").Value = '0.713'
$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '82.73'
$ws.Range("E24").Value = '  -0.41%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '0.0000145'
$ws.Range("E25").Value = '  -3.19%  '

$ws.Range("D26").Value = '12.23'
$ws.Range("E26").Value = '  +1.59%  '

$ws.Range("D27").Value = '10.31'
$ws.Range("E27").Value = '  +2.95%  '

$ws.Range("E28").Value = '  -1.29%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").Value = '3.946.89'
$ws.Range("E30").Value = '  +0.38%  '

$ws.Range("D31").Value = '2.81'
$ws.Range("E31").Value = '  +4.89%  '

$ws.Range("E32").Value = '  +3.51%  '

$ws.Range("D33").Value = '7.36'
$ws.Range("E33").Value = '  +0.77%  '

$ws.Range("D34").Value = '29.16'
$ws.Range("E34").Value = '  +0.65%  '

$ws.Range("D35").Value = '0.179'
$ws.Range("E35").Value = '  +19.31%  '

$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.752.65'
$ws.Range("E37").Value = '  +0.51%  '

$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '9.05'
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D40").Value = '5.90'
$ws.Range("E40").Value = '  +1.21%  '

$ws.Range("D41").Value = '3.30'
$ws.Range("E41").Value = '  -0.74%  '

$ws.Range("D42").Value = '0.963'
$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("E45").Value = '  +6.36%  '

$ws.Range("D46").Value = '45.04'
$ws.Range("E46").Value = '  +5.69%  '

$ws.Range("D47").Value = '158.99'
$ws.Range("E47").Value = '  +3.87%  '

$ws.Range("D48").Value = '48.00'
$ws.Range("E48").Value = '  +2.40%  '

$ws.Range("D49").Value = '0.300'
$ws.Range("E49").Value = '  +0.34%  '

$ws.Range("E50").Value = '  +1.58%  '

$ws.Range("D51").Value = '8.47'
$ws.Range("E51").Value = '  +0.35%  '

# Drop the temporary text-number-format override so the cells end up
# with no explicit style, matching stock (unstyled) data cells.
foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}
